# Automatische test-sync: 2025-08-14 20:55:50
# Append new log entry to the "Logs" sheet and refresh the "Dashboard" count.

$wb = $excel.ActiveWorkbook

# --- 1. Append the new row to the Logs sheet -------------------------------
$logs = $wb.Worksheets.Item("Logs")

$newRow = 16
$logs.Range("A" + $newRow).Value = "Nieuwe bestelling"
$logs.Range("B" + $newRow).Value = "planning@testbedrijf123.nl"
$logs.Range("C" + $newRow).Value = "Wil je 200 stuks M8-bouten bestellen bij onze leverancier?"
$logs.Range("D" + $newRow).Value = "Intern verzoek / Actie voor medewerker"
$logs.Range("E" + $newRow).Value = "Bedankt, we hebben dit doorgestuurd naar inkoop@testbedrijf123.nl."
$logs.Range("F" + $newRow).Value = "2025-08-14 20:54:57"
$logs.Range("G" + $newRow).Value = "Nee"
$logs.Range("H" + $newRow).Value = "Ja"
$logs.Range("I" + $newRow).Value = "Nee"
$logs.Range("J" + $newRow).Value = "Nee"

# --- 2. Extend the conditional-formatting ranges so row 16 is covered ------
$logs.Range("D2:D15").FormatConditions.Item(1).ModifyAppliesToRange($logs.Range("D2:D16"))
$logs.Range("G2:G15").FormatConditions.Item(1).ModifyAppliesToRange($logs.Range("G2:G16"))
$logs.Range("H2:H15").FormatConditions.Item(1).ModifyAppliesToRange($logs.Range("H2:H16"))
$logs.Range("I2:I15").FormatConditions.Item(1).ModifyAppliesToRange($logs.Range("I2:I16"))
$logs.Range("J2:J15").FormatConditions.Item(1).ModifyAppliesToRange($logs.Range("J2:J16"))

# --- 3. Update the Dashboard summary count ----------------------------------
$dashboard = $wb.Worksheets.Item("Dashboard")
$dashboard.Range("B2").Value = 10
